# Auto-generated edit script: updates Leve profit-calc sheets with
# refreshed market-price snapshot values (scheduled runner update).
$wb = $excel.ActiveWorkbook

# ---- ALC ----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H4").Value = 288.75
$ws.Range("I4").Value = 146.5
$ws.Range("J4").Value = 1000
$ws.Range("K4").Value = 146.5
$ws.Range("L4").Value = 1000
$ws.Range("M4").Value = -32.5
$ws.Range("N4").Value = -1228
$ws.Range("H17").Value = 2522.6274
$ws.Range("J17").Value = 2559.82
$ws.Range("L17").Value = 7679.460000000001
$ws.Range("N17").Value = -8015.460000000001
$ws.Range("H28").Value = 256.875
$ws.Range("I28").Value = 329.83334
$ws.Range("J28").Value = 38
$ws.Range("K28").Value = 329.83334
$ws.Range("L28").Value = 38
$ws.Range("M28").Value = 155.16666
$ws.Range("N28").Value = -1008
$ws.Range("H62").Value = 1690.2273
$ws.Range("I62").Value = 1441.9333
$ws.Range("K62").Value = 1441.9333
$ws.Range("M62").Value = -817.9332999999999
$ws.Range("H65").Value = 1690.2273
$ws.Range("I65").Value = 1441.9333
$ws.Range("K65").Value = 7209.666499999999
$ws.Range("M65").Value = -4089.666499999999
$ws.Range("H80").Value = 3600.4167
$ws.Range("I80").Value = 974.2727
$ws.Range("J80").Value = 5822.5386
$ws.Range("K80").Value = 2922.8181
$ws.Range("L80").Value = 17467.6158
$ws.Range("M80").Value = -1924.8181
$ws.Range("N80").Value = -19463.6158
$ws.Range("H83").Value = 3600.4167
$ws.Range("I83").Value = 974.2727
$ws.Range("J83").Value = 5822.5386
$ws.Range("K83").Value = 8768.454299999999
$ws.Range("L83").Value = 52402.8474
$ws.Range("M83").Value = -3776.454299999999
$ws.Range("N83").Value = -62386.8474
$ws.Range("H106").Value = 1739.3529
$ws.Range("I106").Value = 1338.25
$ws.Range("J106").Value = 2702
$ws.Range("K106").Value = 1338.25
$ws.Range("L106").Value = 2702
$ws.Range("M106").Value = -707.25
$ws.Range("N106").Value = -3964
$ws.Range("H107").Value = 1000
$ws.Range("I107").Value = 1000
$ws.Range("K107").Value = 1000
$ws.Range("M107").Value = 920
$ws.Range("H116").Value = 5320
$ws.Range("I116").Value = 998
$ws.Range("J116").Value = 7172.2856
$ws.Range("K116").Value = 998
$ws.Range("L116").Value = 7172.2856
$ws.Range("M116").Value = 2444
$ws.Range("N116").Value = -14056.2856
$ws.Range("H129").Value = 1367.585
$ws.Range("J129").Value = 1401.6078
$ws.Range("L129").Value = 4204.8234
$ws.Range("N129").Value = -14204.8234
$ws.Range("H132").Value = 3580.9565
$ws.Range("I132").Value = 3908.8333
$ws.Range("J132").Value = 2400.6
$ws.Range("K132").Value = 11726.4999
$ws.Range("L132").Value = 7201.799999999999
$ws.Range("M132").Value = -9196.499899999999
$ws.Range("N132").Value = -12261.8
$ws.Range("H137").Value = 92809.55
$ws.Range("I137").Value = 1750.6666
$ws.Range("J137").Value = 202080.2
$ws.Range("K137").Value = 5251.9998
$ws.Range("L137").Value = 606240.6000000001
$ws.Range("M137").Value = -2701.9998
$ws.Range("N137").Value = -611340.6000000001
$ws.Range("H138").Value = 1369.7805
$ws.Range("I138").Value = 540.36664
$ws.Range("J138").Value = 3631.818
$ws.Range("K138").Value = 1621.09992
$ws.Range("L138").Value = 10895.454
$ws.Range("M138").Value = 3518.90008
$ws.Range("N138").Value = -21175.454
$ws.Range("H141").Value = 3782
$ws.Range("I141").Value = 2480
$ws.Range("J141").Value = 5084
$ws.Range("K141").Value = 7440
$ws.Range("L141").Value = 15252
$ws.Range("M141").Value = -2260
$ws.Range("N141").Value = -25612

# ---- ARM ----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 18424.693
$ws.Range("I32").Value = 18152.738
$ws.Range("K32").Value = 18152.738
$ws.Range("M32").Value = -17865.738
$ws.Range("H45").Value = 3761.4849
$ws.Range("I45").Value = 3686.1538
$ws.Range("J45").Value = 3810.45
$ws.Range("K45").Value = 3686.1538
$ws.Range("L45").Value = 3810.45
$ws.Range("M45").Value = -3309.1538
$ws.Range("N45").Value = -4564.45
$ws.Range("H115").Value = 31310.5
$ws.Range("J115").Value = 31310.5
$ws.Range("L115").Value = 31310.5
$ws.Range("N115").Value = -34444.5

# ---- BSM ----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 1395.0892
$ws.Range("I86").Value = 1316.8478
$ws.Range("J86").Value = 1755
$ws.Range("K86").Value = 1316.8478
$ws.Range("L86").Value = 1755
$ws.Range("M86").Value = -193.8478
$ws.Range("N86").Value = -4001
$ws.Range("H89").Value = 1395.0892
$ws.Range("I89").Value = 1316.8478
$ws.Range("J89").Value = 1755
$ws.Range("K89").Value = 6584.239
$ws.Range("L89").Value = 8775
$ws.Range("M89").Value = -968.2389999999996
$ws.Range("N89").Value = -20007
$ws.Range("H107").Value = 744.5
$ws.Range("I107").Value = 659.3333
$ws.Range("K107").Value = 659.3333
$ws.Range("M107").Value = 1260.6667

# ---- CRP ----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 19518.354
$ws.Range("I31").Value = 34757
$ws.Range("J31").Value = 2374.875
$ws.Range("K31").Value = 34757
$ws.Range("L31").Value = 2374.875
$ws.Range("M31").Value = -34462
$ws.Range("N31").Value = -2964.875
$ws.Range("H34").Value = 19518.354
$ws.Range("I34").Value = 34757
$ws.Range("J34").Value = 2374.875
$ws.Range("K34").Value = 34757
$ws.Range("L34").Value = 2374.875
$ws.Range("M34").Value = -34555
$ws.Range("N34").Value = -2778.875
$ws.Range("H58").Value = 19045.5
$ws.Range("I58").Value = 1261.9546
$ws.Range("K58").Value = 1261.9546
$ws.Range("M58").Value = -1058.9546
$ws.Range("H81").Value = 50000
$ws.Range("J81").Value = 50000
$ws.Range("L81").Value = 50000
$ws.Range("N81").Value = -51996
$ws.Range("H84").Value = 50000
$ws.Range("J84").Value = 50000
$ws.Range("L84").Value = 150000
$ws.Range("N84").Value = -159984
$ws.Range("H86").Value = 5214626
$ws.Range("I86").Value = 2022.5
$ws.Range("J86").Value = 11916545
$ws.Range("K86").Value = 2022.5
$ws.Range("L86").Value = 11916545
$ws.Range("M86").Value = -899.5
$ws.Range("N86").Value = -11918791
$ws.Range("H89").Value = 5214626
$ws.Range("I89").Value = 2022.5
$ws.Range("J89").Value = 11916545
$ws.Range("K89").Value = 10112.5
$ws.Range("L89").Value = 59582725
$ws.Range("M89").Value = -4496.5
$ws.Range("N89").Value = -59593957
$ws.Range("H99").Value = 14709702
$ws.Range("I99").Value = 3162.2104
$ws.Range("K99").Value = 3162.2104
$ws.Range("M99").Value = -1664.2104
$ws.Range("H126").Value = 14709702
$ws.Range("I126").Value = 3162.2104
$ws.Range("K126").Value = 9486.6312
$ws.Range("M126").Value = -7016.6312
$ws.Range("H132").Value = 38337.4
$ws.Range("I132").Value = 59283.223
$ws.Range("J132").Value = 6918.6665
$ws.Range("K132").Value = 177849.669
$ws.Range("L132").Value = 20755.9995
$ws.Range("M132").Value = -175319.669
$ws.Range("N132").Value = -25815.9995
$ws.Range("H134").Value = 1152.8889
$ws.Range("I134").Value = 1102
$ws.Range("K134").Value = 3306
$ws.Range("M134").Value = -771
$ws.Range("H136").Value = 19045.5
$ws.Range("I136").Value = 1261.9546
$ws.Range("K136").Value = 3785.8638
$ws.Range("M136").Value = -1235.8638

# ---- CUL ----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H45").Value = 530
$ws.Range("I45").Value = 0
$ws.Range("J45").Value = 530
$ws.Range("K45").Value = 0
$ws.Range("L45").Value = 1590
$ws.Range("N45").Value = -2654
$ws.Range("M45").ClearContents()
$ws.Range("H116").Value = 1200
$ws.Range("I116").Value = 933.3333
$ws.Range("K116").Value = 2799.9999
$ws.Range("M116").Value = 642.0001000000002
$ws.Range("H128").Value = 500000
$ws.Range("I128").Value = 500000
$ws.Range("K128").Value = 1500000
$ws.Range("M128").Value = -1495020
$ws.Range("H131").Value = 735.6
$ws.Range("I131").Value = 700
$ws.Range("J131").Value = 735.9596
$ws.Range("K131").Value = 2100
$ws.Range("L131").Value = 2207.8788
$ws.Range("M131").Value = 2940
$ws.Range("N131").Value = -12287.8788
$ws.Range("H136").Value = 2522.75
$ws.Range("I136").Value = 1713.3334
$ws.Range("J136").Value = 4951
$ws.Range("K136").Value = 5140.0002
$ws.Range("L136").Value = 14853
$ws.Range("M136").Value = -40.0002000000004
$ws.Range("N136").Value = -25053

# ---- GSM ----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 49.63158
$ws.Range("I2").Value = 46.285713
$ws.Range("J2").Value = 59
$ws.Range("K2").Value = 46.285713
$ws.Range("L2").Value = 59
$ws.Range("M2").Value = 66.714287
$ws.Range("N2").Value = -285
$ws.Range("H80").Value = 3664.6875
$ws.Range("I80").Value = 3599.8
$ws.Range("J80").Value = 3694.182
$ws.Range("K80").Value = 3599.8
$ws.Range("L80").Value = 3694.182
$ws.Range("M80").Value = -2601.8
$ws.Range("N80").Value = -5690.182
$ws.Range("H83").Value = 3664.6875
$ws.Range("I83").Value = 3599.8
$ws.Range("J83").Value = 3694.182
$ws.Range("K83").Value = 17999
$ws.Range("L83").Value = 18470.91
$ws.Range("M83").Value = -13007
$ws.Range("N83").Value = -28454.91
$ws.Range("H107").Value = 5495326.5
$ws.Range("J107").Value = 12822245
$ws.Range("L107").Value = 12822245
$ws.Range("N107").Value = -12826085
$ws.Range("H113").Value = 2035.4062
$ws.Range("I113").Value = 1735.2307
$ws.Range("K113").Value = 1735.2307
$ws.Range("M113").Value = 434.7692999999999
$ws.Range("H135").Value = 50585
$ws.Range("I135").Value = 0
$ws.Range("J135").Value = 50585
$ws.Range("K135").Value = 0
$ws.Range("L135").Value = 50585
$ws.Range("N135").Value = -60725
$ws.Range("M135").ClearContents()

# ---- LTW ----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 5991.0557
$ws.Range("I61").Value = 2263
$ws.Range("J61").Value = 8363.454
$ws.Range("K61").Value = 2263
$ws.Range("L61").Value = 8363.454
$ws.Range("M61").Value = -2061
$ws.Range("N61").Value = -8767.454
$ws.Range("H113").Value = 5991.0557
$ws.Range("I113").Value = 2263
$ws.Range("J113").Value = 8363.454
$ws.Range("K113").Value = 2263
$ws.Range("L113").Value = 8363.454
$ws.Range("M113").Value = -93
$ws.Range("N113").Value = -12703.454
$ws.Range("H122").Value = 1637346.1
$ws.Range("I122").Value = 2181161
$ws.Range("J122").Value = 5901.6665
$ws.Range("K122").Value = 6543483
$ws.Range("L122").Value = 17704.9995
$ws.Range("M122").Value = -6541033
$ws.Range("N122").Value = -22604.9995
$ws.Range("H123").Value = 20429
$ws.Range("J123").Value = 20429
$ws.Range("L123").Value = 20429
$ws.Range("N123").Value = -30229
$ws.Range("H136").Value = 46603.547
$ws.Range("I136").Value = 56737.668
$ws.Range("J136").Value = 1000
$ws.Range("K136").Value = 170213.004
$ws.Range("L136").Value = 3000
$ws.Range("M136").Value = -167663.004
$ws.Range("N136").Value = -8100
$ws.Range("H141").Value = 50000
$ws.Range("J141").Value = 50000
$ws.Range("L141").Value = 50000
$ws.Range("N141").Value = -60360

# ---- WVR ----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 1776.5333
$ws.Range("I122").Value = 1843.6666
$ws.Range("K122").Value = 5530.9998
$ws.Range("M122").Value = -3080.9998
$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()
$ws.Range("H136").Value = 37039110
$ws.Range("I136").Value = 55557520
$ws.Range("J136").Value = 2289.3333
$ws.Range("K136").Value = 166672560
$ws.Range("L136").Value = 6867.999899999999
$ws.Range("M136").Value = -166670010
$ws.Range("N136").Value = -11967.9999

